# Apply the updated "cryptos" price-list values described by the commit.
#
# All target cells are plain text cells (inline strings) holding Price /
# Volume(1h) figures (and, for two rows, Coin name + Link), even though many
# of the Price values look numeric (e.g. "109.88"). Assigning such a string
# straight to .Value makes Excel auto-coerce it into a real number, which
# would change the stored cell type/semantics. To avoid that we write every
# value using Excel's leading-apostrophe "text" convention, then immediately
# reset the cell's Style back to "Normal" so no stray quote-prefix / number
# formatting is left behind - the cell ends up holding exactly the plain
# text value with no extra style applied, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '48.779.72'
    'E2' = '  -2.24%  '
    'D3' = '2.616.67'
    'E3' = '  +0.54%  '
    'E4' = '  +0.02%  '
    'D5' = '109.88'
    'E5' = '  +0.44%  '
    'D6' = '322.18'
    'E6' = '  -0.44%  '
    'D7' = '0.524'
    'E7' = '  -1.32%  '
    'E8' = '  +0.00%  '
    'D9' = '0.540'
    'E9' = '  -3.53%  '
    'D10' = '39.44'
    'E10' = '  -3.04%  '
    'D11' = '19.66'
    'E11' = '  -5.15%  '
    'D12' = '0.0808'
    'E12' = '  -1.52%  '
    'E13' = '  +0.28%  '
    'E14' = '  -1.39%  '
    'D15' = '3.023.53'
    'E15' = '  +0.48%  '
    'D16' = '2.623.80'
    'E16' = '  +1.88%  '
    'D17' = '0.860'
    'E17' = '  -0.83%  '
    'D18' = '48.719.61'
    'E18' = '  -2.28%  '
    'D19' = '2.96'
    'E19' = '  -3.75%  '
    'D20' = '12.81'
    'E20' = '  -3.90%  '
    'D21' = '6.67'
    'E21' = '  -1.25%  '
    'E22' = '  -0.86%  '
    'D23' = '268.84'
    'E23' = '  -4.95%  '
    'D24' = '68.58'
    'E24' = '  -5.93%  '
    'E25' = '  -1.03%  '
    'D26' = '26.02'
    'E26' = '  -2.37%  '
    'D27' = '0.999'
    'E27' = '  +0.11%  '
    'D28' = '10.03'
    'E28' = '  +0.68%  '
    'D29' = '2.22'
    'E29' = '  -0.67%  '
    'D30' = '34.79'
    'E30' = '  -2.97%  '
    'E31' = '  -6.99%  '
    'D32' = '49.32'
    'E32' = '  -0.27%  '
    'D33' = '5.47'
    'E33' = '  +0.81%  '
    'E34' = '  -0.38%  '
    'D35' = '0.0796'
    'E35' = '  +0.47%  '
    'D36' = '18.94'
    'E36' = '  -4.30%  '
    'D37' = '5.00'
    'E37' = '  +5.54%  '
    'E38' = '  -1.21%  '
    'E39' = '  +1.73%  '
    'D40' = '126.06'
    'E40' = '  +1.27%  '
    'B41' = 'Stellar'
    'C41' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D41' = '0.111'
    'E41' = '  -1.39%  '
    'B42' = 'EnergySwap'
    'C42' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D42' = '22.25'
    'E42' = '  -2.32%  '
    'D43' = '2.13'
    'E43' = '  -4.42%  '
    'E44' = '  +0.60%  '
    'D45' = '2.057.51'
    'E45' = '  +0.79%  '
    'E46' = '  -3.96%  '
    'E47' = '  +4.94%  '
    'E48' = '  +0.39%  '
    'E49' = '  -3.46%  '
    'D50' = '58.44'
    'E50' = '  +1.72%  '
    'E51' = '  -4.13%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
